# The target diff for this commit ("Fixed POI packaging and upgraded to
# POI 3.15.") does not change any visible/semantic content of the
# document: every hunk re-orders the XML attributes of the same
# elements (e.g. <w:tab w:val="left" w:pos="3119"/> becomes
# <w:tab w:pos="3119" w:val="left"/>, namespace declarations on
# <w:document> get alphabetized, <w:pgSz>/<w:pgMar>/<w:rFonts>/<w:lang>/
# <w:latentStyles>/<w:lsdException>/<w:style>/... attributes get
# alphabetized, etc.) with identical attribute names/values throughout.
#
# That re-ordering is a side effect of the upgraded Apache POI/XMLBeans
# writer used to regenerate the expected-validation fixture at commit
# time (a packaging/library change), not an edit a user performs in
# Word. The Word object model (Find.Execute, Paragraphs, Tables, ...)
# operates above the raw OOXML attribute-serialization layer and has no
# concept of "attribute order" to change, so there is nothing in the
# document's content, formatting or structure for this script to touch.
#
# Applying the *meaning* of the diff therefore means leaving the
# document's content untouched.
$d = $word.ActiveDocument
